$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 95
$ws.Cells.Item($row, 1).Value = "2025-04-29 15:57:03"
$ws.Cells.Item($row, 2).Value = 269
